$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.133.66"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.550.05"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'603.24"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "3.549.17"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "'7.80"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "4.148.05"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'0.0000205"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'29.95"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "3.557.17"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "66.159.48"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  +7.83%  "
$ws.Range("D20").Value = "'6.17"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'14.62"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'428.58"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'0.607"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").Value = "'79.76"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "3.688.73"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "'9.06"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").Value = "'7.81"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "3.543.44"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").Value = "'25.32"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("E35").Value = "  -9.26%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'7.78"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").Value = "'5.52"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "'174.37"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "'0.0844"
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "'5.17"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'0.884"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "'45.90"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "'1.18"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").Value = "'24.68"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").Value = "'2.38"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "'7.08"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "'22.74"
$ws.Range("E51").Value = "  +1.31%  "
